# Updates the NATMI TPM-derived ligand-receptor statistics for Cx3cl1-Itgav
# to reflect a recomputation with the new TPM values ("update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("E2").Value = 3
    $ws.Range("F2").Value = 1
    $ws.Range("G2").Value = 5.682516
    $ws.Range("H2").Value = 17.047548
    $ws.Range("I2").Value = 0.4522589164991918
    $ws.Range("J2").Value = 0.4522589164991919
    $ws.Range("M2").Value = 3.759736666666667
    $ws.Range("N2").Value = 11.27921
    $ws.Range("O2").Value = 0.0683751702595819
    $ws.Range("P2").Value = 0.06837517025958188
    $ws.Range("Q2").Value = 21.36476376412
    $ws.Range("R2").Value = 192.28287387708
    $ws.Range("S2").Value = 0.03092328041704627
    $ws.Range("T2").Value = 0.03092328041704627
    # Row 3
    $ws.Range("E3").Value = 3
    $ws.Range("F3").Value = 1
    $ws.Range("G3").Value = 5.682516
    $ws.Range("H3").Value = 17.047548
    $ws.Range("I3").Value = 0.4522589164991918
    $ws.Range("J3").Value = 0.4522589164991919
    $ws.Range("O3").Value = 0.6514180024294648
    $ws.Range("P3").Value = 0.6514180024294647
    $ws.Range("Q3").Value = 203.544527651836
    $ws.Range("R3").Value = 1831.900748866524
    $ws.Range("S3").Value = 0.2946095999668176
    $ws.Range("T3").Value = 0.2946095999668176
    # Row 4
    $ws.Range("E4").Value = 3
    $ws.Range("F4").Value = 1
    $ws.Range("G4").Value = 5.682516
    $ws.Range("H4").Value = 17.047548
    $ws.Range("I4").Value = 0.4522589164991918
    $ws.Range("J4").Value = 0.4522589164991919
    $ws.Range("O4").Value = 0.2802068273109533
    $ws.Range("P4").Value = 0.2802068273109533
    $ws.Range("Q4").Value = 87.55448283147999
    $ws.Range("R4").Value = 787.99034548332
    $ws.Range("S4").Value = 0.1267260361153279
    $ws.Range("T4").Value = 0.1267260361153279
    # Row 5
    $ws.Range("I5").Value = 0.4336933920535619
    $ws.Range("J5").Value = 0.433693392053562
    $ws.Range("M5").Value = 3.759736666666667
    $ws.Range("N5").Value = 11.27921
    $ws.Range("O5").Value = 0.0683751702595819
    $ws.Range("P5").Value = 0.06837517025958188
    $ws.Range("Q5").Value = 20.48772623215
    $ws.Range("R5").Value = 184.38953608935
    $ws.Range("S5").Value = 0.0296538595221179
    $ws.Range("T5").Value = 0.0296538595221179
    # Row 6
    $ws.Range("I6").Value = 0.4336933920535619
    $ws.Range("J6").Value = 0.433693392053562
    $ws.Range("O6").Value = 0.6514180024294648
    $ws.Range("P6").Value = 0.6514180024294647
    $ws.Range("S6").Value = 0.28251568311839
    $ws.Range("T6").Value = 0.28251568311839
    # Row 7
    $ws.Range("I7").Value = 0.4336933920535619
    $ws.Range("J7").Value = 0.433693392053562
    $ws.Range("O7").Value = 0.2802068273109533
    $ws.Range("P7").Value = 0.2802068273109533
    $ws.Range("S7").Value = 0.121523849413054
    $ws.Range("T7").Value = 0.121523849413054
    # Row 8
    $ws.Range("I8").Value = 0.1140476914472462
    $ws.Range("J8").Value = 0.1140476914472462
    $ws.Range("M8").Value = 3.759736666666667
    $ws.Range("N8").Value = 11.27921
    $ws.Range("O8").Value = 0.0683751702595819
    $ws.Range("P8").Value = 0.06837517025958188
    $ws.Range("Q8").Value = 5.387626195354444
    $ws.Range("R8").Value = 48.48863575818999
    $ws.Range("S8").Value = 0.00779803032041772
    $ws.Range("T8").Value = 0.007798030320417721
    # Row 9
    $ws.Range("I9").Value = 0.1140476914472462
    $ws.Range("J9").Value = 0.1140476914472462
    $ws.Range("O9").Value = 0.6514180024294648
    $ws.Range("P9").Value = 0.6514180024294647
    $ws.Range("S9").Value = 0.07429271934425707
    $ws.Range("T9").Value = 0.07429271934425707
    # Row 10
    $ws.Range("I10").Value = 0.1140476914472462
    $ws.Range("J10").Value = 0.1140476914472462
    $ws.Range("O10").Value = 0.2802068273109533
    $ws.Range("P10").Value = 0.2802068273109533
    $ws.Range("S10").Value = 0.0319569417825714
    $ws.Range("T10").Value = 0.0319569417825714
